# Actualización automática del tracker
# Rellena las columnas "resultado" (G) y "profit" (H) para las filas
# que ya tienen un pronóstico pero aún no tienen resultado cargado.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 58: cuota 2.38 -> Fallo
$ws.Range("G58").Value = "Fallo"
$ws.Range("H58").Value = -1

# Fila 62: cuota 2.63 -> Fallo
$ws.Range("G62").Value = "Fallo"
$ws.Range("H62").Value = -1

# Fila 75: cuota 1.57 -> Acierto (profit = cuota - 1)
$ws.Range("G75").Value = "Acierto"
$ws.Range("H75").Value = 0.57

# Fila 79: cuota 7.5 -> Fallo
$ws.Range("G79").Value = "Fallo"
$ws.Range("H79").Value = -1
